# The "diputados" column (column C) of the metadata sheet is re-curated
# from a dimension to a measure:
#   - row2 (iaest type):     iaest-dimension:diputados -> iaest-measure:diputados
#   - row3 (dim/medida):     dim                        -> medida
#   - row4 (concept type):   skos:Concept               -> xsd:int
#   - row5 (mapping file):   mapping-diputados.xlsx no longer applies -> cell removed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "iaest-measure:diputados"
$ws.Range("C3").Value = "medida"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("C5").Clear()
